# Update Name of Algo
# Apply updated numeric values to result_data_RandomForest sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.650100000000004
$ws.Range("A3").Value = -21.74649999999999
$ws.Range("D3").Value = -7.281799999999996
$ws.Range("D12").Value = -7.337300000000002
$ws.Range("A14").Value = -21.81279999999999
$ws.Range("A16").Value = -21.90249999999999
$ws.Range("B18").Value = 7.014099999999996
$ws.Range("A21").Value = -20.25399999999998
$ws.Range("A23").Value = -20.73259999999998
$ws.Range("B24").Value = 6.393099999999997
$ws.Range("D24").Value = -7.728299999999999
$ws.Range("A25").Value = -21.64269999999999
$ws.Range("B25").Value = 5.4663
$ws.Range("D25").Value = -8.538899999999996
$ws.Range("A26").Value = -21.09369999999997
$ws.Range("B27").Value = 5.732000000000002
$ws.Range("A29").Value = -20.92229999999998
$ws.Range("B30").Value = 6.470499999999999
$ws.Range("B31").Value = 6.1164
$ws.Range("B39").Value = 9.493400000000007
$ws.Range("A40").Value = -19.4657
$ws.Range("D41").Value = -8.277299999999995
$ws.Range("B42").Value = 10.236
$ws.Range("B48").Value = 5.743700000000002
$ws.Range("D50").Value = -8.132800000000001
$ws.Range("B51").Value = 5.800200000000001
$ws.Range("B52").Value = 5.795800000000003
$ws.Range("A53").Value = -22.0298
$ws.Range("D53").Value = -6.089200000000002
$ws.Range("B55").Value = 6.517799999999997
$ws.Range("B56").Value = 5.396199999999999
$ws.Range("D56").Value = -7.8381
$ws.Range("A57").Value = -22.1769
$ws.Range("B57").Value = 4.944299999999997
$ws.Range("D57").Value = -8.306200000000004
$ws.Range("D58").Value = -8.116099999999999
$ws.Range("A59").Value = -22.4696
$ws.Range("B60").Value = 5.667699999999999
$ws.Range("D61").Value = -7.904399999999997
$ws.Range("D63").Value = -7.927400000000003
$ws.Range("D64").Value = -8.165299999999995
$ws.Range("A65").Value = -21.78279999999998
$ws.Range("A69").Value = -21.59009999999999
$ws.Range("D70").Value = -7.238899999999997
$ws.Range("D72").Value = -7.326400000000008
$ws.Range("B73").Value = 8.589799999999999
$ws.Range("B74").Value = 9.762199999999993
$ws.Range("A79").Value = -20.4516
$ws.Range("A83").Value = -21.7993
$ws.Range("D86").Value = -7.7933
$ws.Range("B89").Value = 4.85489999999999
$ws.Range("D89").Value = -5.648500000000001
$ws.Range("B90").Value = 5.935000000000002
$ws.Range("A91").Value = -21.40080000000002
$ws.Range("B92").Value = 4.944099999999992
$ws.Range("A93").Value = -21.13979999999999
$ws.Range("D98").Value = -8.510500000000002
$ws.Range("A100").Value = -21.71919999999999
$ws.Range("D100").Value = -8.1897
$ws.Range("D102").Value = -7.664299999999996
